$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-07 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-08 Saturday", 2) | Out-Null
$d.Content.Find.Execute("251×7=", $true, $false, $false, $false, $false, $true, 1, $false, "434×9=", 2) | Out-Null
$d.Content.Find.Execute("532×7=", $true, $false, $false, $false, $false, $true, 1, $false, "714×5=", 2) | Out-Null
$d.Content.Find.Execute("966×3=", $true, $false, $false, $false, $false, $true, 1, $false, "439×3=", 2) | Out-Null
$d.Content.Find.Execute("768×9=", $true, $false, $false, $false, $false, $true, 1, $false, "449×7=", 2) | Out-Null
$d.Content.Find.Execute("692×7=", $true, $false, $false, $false, $false, $true, 1, $false, "673×6=", 2) | Out-Null
$d.Content.Find.Execute("398×9=", $true, $false, $false, $false, $false, $true, 1, $false, "435×2=", 2) | Out-Null
$d.Content.Find.Execute("554×4=", $true, $false, $false, $false, $false, $true, 1, $false, "603×4=", 2) | Out-Null
$d.Content.Find.Execute("514×6=", $true, $false, $false, $false, $false, $true, 1, $false, "385×4=", 2) | Out-Null
$d.Content.Find.Execute("868×3=", $true, $false, $false, $false, $false, $true, 1, $false, "915×7=", 2) | Out-Null
$d.Content.Find.Execute("820×5=", $true, $false, $false, $false, $false, $true, 1, $false, "495×4=", 2) | Out-Null
$d.Content.Find.Execute("553×4=", $true, $false, $false, $false, $false, $true, 1, $false, "263×4=", 2) | Out-Null
$d.Content.Find.Execute("412×3=", $true, $false, $false, $false, $false, $true, 1, $false, "510×5=", 2) | Out-Null
$d.Content.Find.Execute("215×8=", $true, $false, $false, $false, $false, $true, 1, $false, "758×4=", 2) | Out-Null
$d.Content.Find.Execute("529×7=", $true, $false, $false, $false, $false, $true, 1, $false, "267×4=", 2) | Out-Null
$d.Content.Find.Execute("901×3=", $true, $false, $false, $false, $false, $true, 1, $false, "129×4=", 2) | Out-Null
$d.Content.Find.Execute("298×5=", $true, $false, $false, $false, $false, $true, 1, $false, "266×4=", 2) | Out-Null
$d.Content.Find.Execute("963×5=", $true, $false, $false, $false, $false, $true, 1, $false, "572×2=", 2) | Out-Null
$d.Content.Find.Execute("370×7=", $true, $false, $false, $false, $false, $true, 1, $false, "179×2=", 2) | Out-Null
$d.Content.Find.Execute("836×4=", $true, $false, $false, $false, $false, $true, 1, $false, "434×5=", 2) | Out-Null
$d.Content.Find.Execute("437×8=", $true, $false, $false, $false, $false, $true, 1, $false, "765×8=", 2) | Out-Null
$d.Content.Find.Execute("914×9=", $true, $false, $false, $false, $false, $true, 1, $false, "197×8=", 2) | Out-Null
$d.Content.Find.Execute("166×5=", $true, $false, $false, $false, $false, $true, 1, $false, "835×5=", 2) | Out-Null
$d.Content.Find.Execute("184×8=", $true, $false, $false, $false, $false, $true, 1, $false, "195×4=", 2) | Out-Null
$d.Content.Find.Execute("228×5=", $true, $false, $false, $false, $false, $true, 1, $false, "553×2=", 2) | Out-Null
$d.Content.Find.Execute("892×9=", $true, $false, $false, $false, $false, $true, 1, $false, "479×4=", 2) | Out-Null
